# Auto-generated edit script: update cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.371.35"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.600.17"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.333"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").Value = "3.061.67"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "59.295.79"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000133"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.564.12"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.408"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.164"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("D28").Value = "0.0₃0742"
$ws.Range("E28").Value = "  +2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.15%  "
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.830"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.823"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "272.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0954"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.19%  "
$ws.Range("D47").Value = "1.942.80"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0222"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("E51").Value = "  +0.30%  "
